$d = $word.ActiveDocument
$d.Content.Find.Execute("Ben Bar", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Ben Barrrr", 2)
